$wb = $excel.ActiveWorkbook

# Sheet ALC, row 116 (@@ -6454,25 +6454,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 1809.5454  # H116: 1955.6666 -> 1809.5454
$ws.Cells.Item(116, 9).Value = 1683.1666  # I116: 1965 -> 1683.1666
$ws.Cells.Item(116, 10).Value = 1961.2  # J116: 1951 -> 1961.2
$ws.Cells.Item(116, 11).Value = 1683.1666  # K116: 1965 -> 1683.1666
$ws.Cells.Item(116, 12).Value = 1961.2  # L116: 1951 -> 1961.2
$ws.Cells.Item(116, 13).Value = 1758.8334  # M116: 1477 -> 1758.8334
$ws.Cells.Item(116, 14).Value = -8845.200000000001  # N116: -8835 -> -8845.200000000001

# Sheet ALC, row 129 (@@ -7091,25 +7091,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 16268.877  # H129: 40819.84 -> 16268.877
$ws.Cells.Item(129, 10).Value = 20989.5  # J129: 101249.4 -> 20989.5
$ws.Cells.Item(129, 12).Value = 62968.5  # L129: 303748.2 -> 62968.5
$ws.Cells.Item(129, 14).Value = -72968.5  # N129: -313748.2 -> -72968.5

# Sheet ALC, row 132 (@@ -7241,25 +7241,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 4609905.5  # H132: 4763559 -> 4609905.5
$ws.Cells.Item(132, 9).Value = 4927796  # I132: 4927802.5 -> 4927796
$ws.Cells.Item(132, 10).Value = 498.5  # J132: 500 -> 498.5
$ws.Cells.Item(132, 11).Value = 14783388  # K132: 14783407.5 -> 14783388
$ws.Cells.Item(132, 12).Value = 1495.5  # L132: 1500 -> 1495.5
$ws.Cells.Item(132, 13).Value = -14780858  # M132: -14780877.5 -> -14780858
$ws.Cells.Item(132, 14).Value = -6555.5  # N132: -6560 -> -6555.5

# Sheet ALC, row 135 (@@ -7391,25 +7391,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 767.8077  # H135: 721.4706 -> 767.8077
$ws.Cells.Item(135, 9).Value = 782.52  # I135: 719.0909 -> 782.52
$ws.Cells.Item(135, 10).Value = 400  # J135: 800 -> 400
$ws.Cells.Item(135, 11).Value = 7042.68  # K135: 6471.8181 -> 7042.68
$ws.Cells.Item(135, 12).Value = 3600  # L135: 7200 -> 3600
$ws.Cells.Item(135, 13).Value = -4507.68  # M135: -3936.8181 -> -4507.68
$ws.Cells.Item(135, 14).Value = -8670  # N135: -12270 -> -8670

# Sheet ALC, row 141 (@@ -7688,25 +7688,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 2862.9656  # H141: 2692.9443 -> 2862.9656
$ws.Cells.Item(141, 9).Value = 1357.2222  # I141: 1265.909 -> 1357.2222
$ws.Cells.Item(141, 10).Value = 5326.909  # J141: 4935.4287 -> 5326.909
$ws.Cells.Item(141, 11).Value = 4071.6666  # K141: 3797.727 -> 4071.6666
$ws.Cells.Item(141, 12).Value = 15980.727  # L141: 14806.2861 -> 15980.727
$ws.Cells.Item(141, 13).Value = 1108.3334  # M141: 1382.273 -> 1108.3334
$ws.Cells.Item(141, 14).Value = -26340.727  # N141: -25166.2861 -> -26340.727

# Sheet ARM, row 2 (@@ -7834,25 +7834,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1327.125  # H2: 1333.375 -> 1327.125
$ws.Cells.Item(2, 9).Value = 1103.7778  # I2: 1043.4 -> 1103.7778
$ws.Cells.Item(2, 10).Value = 1614.2858  # J2: 1816.6666 -> 1614.2858
$ws.Cells.Item(2, 11).Value = 1103.7778  # K2: 1043.4 -> 1103.7778
$ws.Cells.Item(2, 12).Value = 1614.2858  # L2: 1816.6666 -> 1614.2858
$ws.Cells.Item(2, 13).Value = -990.7778000000001  # M2: -930.4000000000001 -> -990.7778000000001
$ws.Cells.Item(2, 14).Value = -1840.2858  # N2: -2042.6666 -> -1840.2858

# Sheet ARM, row 5 (@@ -7984,25 +7984,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 682.7143  # H5: 20833996 -> 682.7143
$ws.Cells.Item(5, 9).Value = 500  # I5: 41667104 -> 500
$ws.Cells.Item(5, 10).Value = 755.8  # J5: 887.5 -> 755.8
$ws.Cells.Item(5, 11).Value = 500  # K5: 41667104 -> 500
$ws.Cells.Item(5, 12).Value = 755.8  # L5: 887.5 -> 755.8
$ws.Cells.Item(5, 13).Value = -388  # M5: -41666992 -> -388
$ws.Cells.Item(5, 14).Value = -979.8  # N5: -1111.5 -> -979.8

# Sheet ARM, row 32 (@@ -9274,22 +9274,22 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3626.0793  # H32: 4305 -> 3626.0793
$ws.Cells.Item(32, 9).Value = 3023.5615  # I32: 3646.9565 -> 3023.5615
$ws.Cells.Item(32, 11).Value = 3023.5615  # K32: 3646.9565 -> 3023.5615
$ws.Cells.Item(32, 13).Value = -2736.5615  # M32: -3359.9565 -> -2736.5615

# Sheet ARM, row 116 (@@ -13384,25 +13384,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 1327.125  # H116: 1333.375 -> 1327.125
$ws.Cells.Item(116, 9).Value = 1103.7778  # I116: 1043.4 -> 1103.7778
$ws.Cells.Item(116, 10).Value = 1614.2858  # J116: 1816.6666 -> 1614.2858
$ws.Cells.Item(116, 11).Value = 1103.7778  # K116: 1043.4 -> 1103.7778
$ws.Cells.Item(116, 12).Value = 1614.2858  # L116: 1816.6666 -> 1614.2858
$ws.Cells.Item(116, 13).Value = 1190.2222  # M116: 1250.6 -> 1190.2222
$ws.Cells.Item(116, 14).Value = -6202.2858  # N116: -6404.6666 -> -6202.2858

# Sheet BSM, row 3 (@@ -14795,25 +14795,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1327.125  # H3: 1333.375 -> 1327.125
$ws.Cells.Item(3, 9).Value = 1103.7778  # I3: 1043.4 -> 1103.7778
$ws.Cells.Item(3, 10).Value = 1614.2858  # J3: 1816.6666 -> 1614.2858
$ws.Cells.Item(3, 11).Value = 1103.7778  # K3: 1043.4 -> 1103.7778
$ws.Cells.Item(3, 12).Value = 1614.2858  # L3: 1816.6666 -> 1614.2858
$ws.Cells.Item(3, 13).Value = -989.7778000000001  # M3: -929.4000000000001 -> -989.7778000000001
$ws.Cells.Item(3, 14).Value = -1842.2858  # N3: -2044.6666 -> -1842.2858

# Sheet BSM, row 4 (@@ -14847,25 +14847,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 682.7143  # H4: 20833996 -> 682.7143
$ws.Cells.Item(4, 9).Value = 500  # I4: 41667104 -> 500
$ws.Cells.Item(4, 10).Value = 755.8  # J4: 887.5 -> 755.8
$ws.Cells.Item(4, 11).Value = 500  # K4: 41667104 -> 500
$ws.Cells.Item(4, 12).Value = 755.8  # L4: 887.5 -> 755.8
$ws.Cells.Item(4, 13).Value = -385  # M4: -41666989 -> -385
$ws.Cells.Item(4, 14).Value = -985.8  # N4: -1117.5 -> -985.8

# Sheet BSM, row 20 (@@ -15628,25 +15628,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3429.9546  # H20: 3447.318 -> 3429.9546
$ws.Cells.Item(20, 9).Value = 3883.2307  # I20: 4493.091 -> 3883.2307
$ws.Cells.Item(20, 10).Value = 2775.2222  # J20: 2401.5454 -> 2775.2222
$ws.Cells.Item(20, 11).Value = 3883.2307  # K20: 4493.091 -> 3883.2307
$ws.Cells.Item(20, 12).Value = 2775.2222  # L20: 2401.5454 -> 2775.2222
$ws.Cells.Item(20, 13).Value = -3636.2307  # M20: -4246.091 -> -3636.2307
$ws.Cells.Item(20, 14).Value = -3269.2222  # N20: -2895.5454 -> -3269.2222

# Sheet BSM, row 22 (@@ -15729,22 +15729,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 353.54544  # H22: 444.44446 -> 353.54544
$ws.Cells.Item(22, 9).Value = 376.7647  # I22: 444.44446 -> 376.7647
$ws.Cells.Item(22, 10).Value = 274.6  # J22: 0 -> 274.6
$ws.Cells.Item(22, 11).Value = 376.7647  # K22: 444.44446 -> 376.7647
$ws.Cells.Item(22, 12).Value = 274.6  # L22: 0 -> 274.6
$ws.Cells.Item(22, 13).Value = -203.7647  # M22: -271.44446 -> -203.7647
$ws.Cells.Item(22, 14).Value = -620.6  # N22: None -> -620.6

# Sheet BSM, row 86 (@@ -18835,25 +18838,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1854.4286  # H86: 2853.8572 -> 1854.4286
$ws.Cells.Item(86, 9).Value = 1612.4445  # I86: 2162.8333 -> 1612.4445
$ws.Cells.Item(86, 10).Value = 2290  # J86: 7000 -> 2290
$ws.Cells.Item(86, 11).Value = 1612.4445  # K86: 2162.8333 -> 1612.4445
$ws.Cells.Item(86, 12).Value = 2290  # L86: 7000 -> 2290
$ws.Cells.Item(86, 13).Value = -489.4445000000001  # M86: -1039.8333 -> -489.4445000000001
$ws.Cells.Item(86, 14).Value = -4536  # N86: -9246 -> -4536

# Sheet BSM, row 89 (@@ -18979,25 +18982,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 1854.4286  # H89: 2853.8572 -> 1854.4286
$ws.Cells.Item(89, 9).Value = 1612.4445  # I89: 2162.8333 -> 1612.4445
$ws.Cells.Item(89, 10).Value = 2290  # J89: 7000 -> 2290
$ws.Cells.Item(89, 11).Value = 8062.2225  # K89: 10814.1665 -> 8062.2225
$ws.Cells.Item(89, 12).Value = 11450  # L89: 35000 -> 11450
$ws.Cells.Item(89, 13).Value = -2446.2225  # M89: -5198.166499999999 -> -2446.2225
$ws.Cells.Item(89, 14).Value = -22682  # N89: -46232 -> -22682

# Sheet BSM, row 105 (@@ -19763,22 +19766,22 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 4167.778  # H105: 4001 -> 4167.778
$ws.Cells.Item(105, 9).Value = 4430  # I105: 4188.75 -> 4430
$ws.Cells.Item(105, 11).Value = 4430  # K105: 4188.75 -> 4430
$ws.Cells.Item(105, 13).Value = -2683  # M105: -2441.75 -> -2683

# Sheet BSM, row 134 (@@ -21169,25 +21172,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 25332.186  # H134: 28377.764 -> 25332.186
$ws.Cells.Item(134, 9).Value = 41964  # I134: 37466.605 -> 41964
$ws.Cells.Item(134, 10).Value = 2232.4443  # J134: 2929 -> 2232.4443
$ws.Cells.Item(134, 11).Value = 125892  # K134: 112399.815 -> 125892
$ws.Cells.Item(134, 12).Value = 6697.3329  # L134: 8787 -> 6697.3329
$ws.Cells.Item(134, 13).Value = -123357  # M134: -109864.815 -> -123357
$ws.Cells.Item(134, 14).Value = -11767.3329  # N134: -13857 -> -11767.3329

# Sheet CRP, row 62 (@@ -24652,25 +24655,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 66669370  # H62: 27780156 -> 66669370
$ws.Cells.Item(62, 9).Value = 2500  # I62: 2087.5 -> 2500
$ws.Cells.Item(62, 10).Value = 83336090  # J62: 55558224 -> 83336090
$ws.Cells.Item(62, 11).Value = 2500  # K62: 2087.5 -> 2500
$ws.Cells.Item(62, 12).Value = 83336090  # L62: 55558224 -> 83336090
$ws.Cells.Item(62, 13).Value = -1876  # M62: -1463.5 -> -1876
$ws.Cells.Item(62, 14).Value = -83337338  # N62: -55559472 -> -83337338

# Sheet CRP, row 65 (@@ -24796,25 +24799,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 66669370  # H65: 27780156 -> 66669370
$ws.Cells.Item(65, 9).Value = 2500  # I65: 2087.5 -> 2500
$ws.Cells.Item(65, 10).Value = 83336090  # J65: 55558224 -> 83336090
$ws.Cells.Item(65, 11).Value = 12500  # K65: 10437.5 -> 12500
$ws.Cells.Item(65, 12).Value = 416680450  # L65: 277791120 -> 416680450
$ws.Cells.Item(65, 13).Value = -9380  # M65: -7317.5 -> -9380
$ws.Cells.Item(65, 14).Value = -416686690  # N65: -277797360 -> -416686690

# Sheet CRP, row 94 (@@ -26202,22 +26205,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 399  # H94: 398 -> 399
$ws.Cells.Item(94, 9).Value = 399  # I94: 398 -> 399
$ws.Cells.Item(94, 11).Value = 399  # K94: 398 -> 399
$ws.Cells.Item(94, 13).Value = 52  # M94: 53 -> 52

# Sheet CUL, row 41 (@@ -30619,22 +30622,19 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(41, 8).Value = 0  # H41: 950 -> 0
$ws.Cells.Item(41, 10).Value = 0  # J41: 950 -> 0
$ws.Cells.Item(41, 12).Value = 0  # L41: 2850 -> 0
$ws.Cells.Item(41, 14).ClearContents()  # N41: -3526 -> (removed)

# Sheet CUL, row 69 (@@ -32018,22 +32018,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 1538  # H69: 1412.5 -> 1538
$ws.Cells.Item(69, 9).Value = 808  # I69: 0 -> 808
$ws.Cells.Item(69, 10).Value = 1684  # J69: 1412.5 -> 1684
$ws.Cells.Item(69, 11).Value = 2424  # K69: 0 -> 2424
$ws.Cells.Item(69, 12).Value = 5052  # L69: 4237.5 -> 5052
$ws.Cells.Item(69, 13).Value = -1613  # M69: None -> -1613
$ws.Cells.Item(69, 14).Value = -6674  # N69: -5859.5 -> -6674

# Sheet CUL, row 72 (@@ -32171,22 +32174,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(72, 8).Value = 1538  # H72: 1412.5 -> 1538
$ws.Cells.Item(72, 9).Value = 808  # I72: 0 -> 808
$ws.Cells.Item(72, 10).Value = 1684  # J72: 1412.5 -> 1684
$ws.Cells.Item(72, 11).Value = 7272  # K72: 0 -> 7272
$ws.Cells.Item(72, 12).Value = 15156  # L72: 12712.5 -> 15156
$ws.Cells.Item(72, 13).Value = -3216  # M72: None -> -3216
$ws.Cells.Item(72, 14).Value = -23268  # N72: -20824.5 -> -23268

# Sheet CUL, row 93 (@@ -33230,22 +33236,22 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(93, 8).Value = 5000  # H93: 3100 -> 5000
$ws.Cells.Item(93, 10).Value = 5000  # J93: 3100 -> 5000
$ws.Cells.Item(93, 12).Value = 15000  # L93: 9300 -> 15000
$ws.Cells.Item(93, 14).Value = -18744  # N93: -13044 -> -18744

# Sheet GSM, row 126 (@@ -41828,25 +41834,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2028.6428  # H126: 1639.1428 -> 2028.6428
$ws.Cells.Item(126, 9).Value = 2067.5833  # I126: 1672.5555 -> 2067.5833
$ws.Cells.Item(126, 10).Value = 1795  # J126: 1579 -> 1795
$ws.Cells.Item(126, 11).Value = 6202.749899999999  # K126: 5017.666499999999 -> 6202.749899999999
$ws.Cells.Item(126, 12).Value = 5385  # L126: 4737 -> 5385
$ws.Cells.Item(126, 13).Value = -3732.749899999999  # M126: -2547.666499999999 -> -3732.749899999999
$ws.Cells.Item(126, 14).Value = -10325  # N126: -9677 -> -10325

# Sheet LTW, row 61 (@@ -45588,25 +45594,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2170  # H61: 2181.7273 -> 2170
$ws.Cells.Item(61, 9).Value = 1925  # I61: 1833.1666 -> 1925
$ws.Cells.Item(61, 10).Value = 2333.3333  # J61: 2600 -> 2333.3333
$ws.Cells.Item(61, 11).Value = 1925  # K61: 1833.1666 -> 1925
$ws.Cells.Item(61, 12).Value = 2333.3333  # L61: 2600 -> 2333.3333
$ws.Cells.Item(61, 13).Value = -1723  # M61: -1631.1666 -> -1723
$ws.Cells.Item(61, 14).Value = -2737.3333  # N61: -3004 -> -2737.3333

# Sheet LTW, row 93 (@@ -47153,25 +47159,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1331.0667  # H93: 1975.375 -> 1331.0667
$ws.Cells.Item(93, 9).Value = 1160.5454  # I93: 2201 -> 1160.5454
$ws.Cells.Item(93, 10).Value = 1800  # J93: 1840 -> 1800
$ws.Cells.Item(93, 11).Value = 1160.5454  # K93: 2201 -> 1160.5454
$ws.Cells.Item(93, 12).Value = 1800  # L93: 1840 -> 1800
$ws.Cells.Item(93, 13).Value = 87.45460000000003  # M93: -953 -> 87.45460000000003
$ws.Cells.Item(93, 14).Value = -4296  # N93: -4336 -> -4296

# Sheet LTW, row 110 (@@ -47977,22 +47983,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(110, 8).Value = 20644  # H110: 19644 -> 20644
$ws.Cells.Item(110, 10).Value = 20644  # J110: 19644 -> 20644
$ws.Cells.Item(110, 12).Value = 20644  # L110: 19644 -> 20644
$ws.Cells.Item(110, 14).Value = -28824  # N110: -27824 -> -28824

# Sheet LTW, row 113 (@@ -48121,25 +48127,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 2170  # H113: 2181.7273 -> 2170
$ws.Cells.Item(113, 9).Value = 1925  # I113: 1833.1666 -> 1925
$ws.Cells.Item(113, 10).Value = 2333.3333  # J113: 2600 -> 2333.3333
$ws.Cells.Item(113, 11).Value = 1925  # K113: 1833.1666 -> 1925
$ws.Cells.Item(113, 12).Value = 2333.3333  # L113: 2600 -> 2333.3333
$ws.Cells.Item(113, 13).Value = 245  # M113: 336.8334 -> 245
$ws.Cells.Item(113, 14).Value = -6673.3333  # N113: -6940 -> -6673.3333

# Sheet LTW, row 132 (@@ -49034,25 +49040,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2403.5557  # H132: 2417.1794 -> 2403.5557
$ws.Cells.Item(132, 9).Value = 2267.276  # I132: 2316.8572 -> 2267.276
$ws.Cells.Item(132, 10).Value = 2968.1428  # J132: 2672.5454 -> 2968.1428
$ws.Cells.Item(132, 11).Value = 6801.828  # K132: 6950.571599999999 -> 6801.828
$ws.Cells.Item(132, 12).Value = 8904.428400000001  # L132: 8017.6362 -> 8904.428400000001
$ws.Cells.Item(132, 13).Value = -4271.828  # M132: -4420.571599999999 -> -4271.828
$ws.Cells.Item(132, 14).Value = -13964.4284  # N132: -13077.6362 -> -13964.4284

# Sheet WVR, row 113 (@@ -55027,25 +55033,25 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 714.75  # H113: 620.4375 -> 714.75
$ws.Cells.Item(113, 9).Value = 799.8570999999999  # I113: 654.9 -> 799.8570999999999
$ws.Cells.Item(113, 10).Value = 595.6  # J113: 563 -> 595.6
$ws.Cells.Item(113, 11).Value = 2399.5713  # K113: 1964.7 -> 2399.5713
$ws.Cells.Item(113, 12).Value = 1786.8  # L113: 1689 -> 1786.8
$ws.Cells.Item(113, 13).Value = -229.5712999999996  # M113: 205.3000000000002 -> -229.5712999999996
$ws.Cells.Item(113, 14).Value = -6126.8  # N113: -6029 -> -6126.8

# Sheet WVR, row 132 (@@ -55958,25 +55964,25 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2569.8823  # H132: 2325.4736 -> 2569.8823
$ws.Cells.Item(132, 9).Value = 2049.3572  # I132: 2380.6365 -> 2049.3572
$ws.Cells.Item(132, 10).Value = 4999  # J132: 2249.625 -> 4999
$ws.Cells.Item(132, 11).Value = 6148.071599999999  # K132: 7141.9095 -> 6148.071599999999
$ws.Cells.Item(132, 12).Value = 14997  # L132: 6748.875 -> 14997
$ws.Cells.Item(132, 13).Value = -3618.071599999999  # M132: -4611.9095 -> -3618.071599999999
$ws.Cells.Item(132, 14).Value = -20057  # N132: -11808.875 -> -20057
